# Remove photo from PPTX.
#
# Slide 9 ("Dr. Jo Guldi's use of data science") has a picture of
# Dr. Jo Guldi (cNvPr id="1026" name="Picture 2",
# descr="http://blog.smu.edu/forum/files/2017/04/jo-guldi-450.jpg")
# that should be deleted, leaving only the title and content
# placeholder shapes on the slide.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(9)

for ($i = $s.Shapes.Count; $i -ge 1; $i--) {
    $sh = $s.Shapes.Item($i)
    if ($sh.Id -eq 1026 -and $sh.Name -eq "Picture 2") {
        $sh.Delete()
    }
}
